$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Shift rows 11-14 down to 12-15, one row at a time from the bottom up, so
# that the source/destination ranges used in each Copy/PasteSpecial never
# span more than a single row (multi-row overlapping copy/paste is mishandled
# by the engine). Blank source cells don't clear the destination on paste,
# so any column-B cell that was blank in the source is cleared explicitly
# afterwards.
#
# Row 15 is brand new (beyond the sheet's old A1:B14 dimension); a value
# paste alone leaves such cells on the default style, so its formatting is
# re-applied with a second, format-only paste from the same source.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4104)
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

$ws.Range("A13:B13").Copy()
$ws.Range("A14:B14").PasteSpecial(-4104)
$ws.Cells.Item(14, 2).ClearContents()

$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4104)
$ws.Cells.Item(13, 2).ClearContents()

$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4104)

# Fill the newly freed row 11 with the "Jurisdiction" entry, copying the
# formatting of an existing body row first so no new cell style gets created.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4104)
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "iso:code:3166:FR"

# Update the Version value (row 3, column B)
$ws.Cells.Item(3, 2).Value = "0.2.0"

# Update the Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2023-10-20T08:59:58+00:00"
